$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New week block: "Jeudi 13 Février 2020" / "Vendredi 14 Février 2020"
# mirrors the existing blocks (rows 1-6, 12-17) at rows 26-31.
# ---------------------------------------------------------------------------

# Row 26 : the two date headers (style copied from the existing "Avertissement"
# header cell B1 so we reuse the same red, border-less style). The Vendredi
# header is entered before the Jeudi one to match the shared-string order.
$ws.Range("F1").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("F26").Value = "                    Vendredi 14 Février 2020"

$ws.Range("B1").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("B26").Value = "                    Jeudi 13 Février 2020"

# Row 27 : column headers (Heure / Objectif / Commentaire) x2, bordered style
# copied from the existing header row (row 2).
$ws.Range("A2:C2").Copy()
$ws.Range("A27:C27").PasteSpecial(-4122)
$ws.Range("E2:G2").Copy()
$ws.Range("E27:G27").PasteSpecial(-4122)

$ws.Range("A27").Value = "Heure"
$ws.Range("B27").Value = "Objectif"
$ws.Range("C27").Value = "Commentaire"
$ws.Range("E27").Value = "Heure"
$ws.Range("F27").Value = "Objectif"
$ws.Range("G27").Value = "Commentaire"

# Rows 28-31 : data rows, bordered style copied from the matching existing
# data rows (3-6) so the border formatting matches exactly.
$ws.Range("A3:C3").Copy()
$ws.Range("A28:C28").PasteSpecial(-4122)
$ws.Range("E3:G3").Copy()
$ws.Range("E28:G28").PasteSpecial(-4122)

$ws.Range("A4:C4").Copy()
$ws.Range("A29:C29").PasteSpecial(-4122)
$ws.Range("E4:G4").Copy()
$ws.Range("E29:G29").PasteSpecial(-4122)

$ws.Range("A5:C5").Copy()
$ws.Range("A30:C30").PasteSpecial(-4122)
$ws.Range("E5:G5").Copy()
$ws.Range("E30:G30").PasteSpecial(-4122)

$ws.Range("A6:C6").Copy()
$ws.Range("A31:C31").PasteSpecial(-4122)
$ws.Range("E6:G6").Copy()
$ws.Range("E31:G31").PasteSpecial(-4122)

# Hour labels (left block)
$ws.Range("A28").Value = "9h30-10h30"
$ws.Range("A29").Value = "10h30-12h30"
$ws.Range("A30").Value = "13h30-15h"
$ws.Range("A31").Value = "15h-16h30"

# Objectif text (left block)
$ws.Range("B28").Value = "Revoir les objectifs"
$ws.Range("B29").Value = "Etablir une section avis dans la page contact"
$ws.Range("B30").Value = "Commencer les modifications du site"
$ws.Range("B31").Value = "Refaire un document de base de données"

# Commentaire (left block) - all done
$ws.Range("C28").Value = "Fait"
$ws.Range("C29").Value = "Fait"
$ws.Range("C30").Value = "Fait"
$ws.Range("C31").Value = "Fait"

# Hour labels (right block)
$ws.Range("E28").Value = "9h30-10h30"
$ws.Range("E29").Value = "10h30-12h30"
$ws.Range("E30").Value = "13h30-15h"
$ws.Range("E31").Value = "15h-16h30"

# F28:F31 stay empty (no Objectif entered yet for the Friday block) but keep
# the bordered formatting already pasted above.

# Commentaire (right block) - all done
$ws.Range("G28").Value = "Fait"
$ws.Range("G29").Value = "Fait"
$ws.Range("G30").Value = "Fait"
$ws.Range("G31").Value = "Fait"

# ---------------------------------------------------------------------------
# View state: selection moves to D20 (matches the author's last click while
# editing the new block).
# ---------------------------------------------------------------------------
$ws.Range("D20").Select()
